$sheetSizes = @{
    "Bags_df" = @('N/A', 'N/A', 'N/A', 'N/A', 'N/A', 'N/A')
    "Skate_df" = @('129,139,149', '54MM', '8 1/2', '8 1/8')
    "Shirts_df" = @('Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge')
    "Pants_df" = @('Small,Medium,Large,XLarge', '30,32,34,36,38', '30,32,34,36,38', '30,32,34,36,38', 'Small,Medium,Large,XLarge,XXL', '30,32,34,36,38', 'Small,Medium,Large,XLarge,XXL', '30,32,34,36,38', '30,32,34,36,38', '30,32,34,36,38', '30,32,34,36,38')
    "Shorts_df" = @('Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge,XXL')
    "Tops_Sweaters_df" = @('Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL')
    "T_Shirts_df" = @('Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL')
    "Jackets_df" = @('Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge', 'Small,Medium,Large,XLarge,XXL')
    "Sweatshirts_df" = @('Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL')
    "Hats_df" = @('N/A', 'N/A', 'N/A', 'S/M,M/L', 'N/A', 'N/A')
    "Accessories_df" = @('N/A', 'N/A', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'Small,Medium,Large,XLarge,XXL', 'N/A')
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $sheetSizes.ContainsKey($name)) {
        continue
    }
    $sizes = $sheetSizes[$name]

    # Header cell: label + formatting copied from the existing "Colors" header (column K)
    $headerCol = $ws.UsedRange.Columns.Count + 1
    $headerCell = $ws.Cells.Item(1, $headerCol)
    $colorsHeaderCell = $ws.Cells.Item(1, $headerCol - 1)
    $colorsHeaderCell.Copy()
    $headerCell.PasteSpecial(-4122)
    $headerCell.Value = "Sizes"

    for ($i = 0; $i -lt $sizes.Length; $i++) {
        $row = 2 + $i
        $ws.Cells.Item($row, $headerCol).Value = $sizes[$i]
    }
}

$excel.CutCopyMode = $false
